$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208, shifting the existing rows 208-268 down to 209-269.
$ws.Rows("208:208").Insert()

# Populate the newly inserted row 208 with a new price-report entry
# (same market/category/quality/prices as the row that used to occupy 208,
# but dated 2022-10-07 / serial 44841).
$ws.Range("A208").Value = 4
$ws.Range("B208").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C208").Value = "Los Lagos"
$ws.Range("D208").Value = 44841
$ws.Range("E208").Value = 10
$ws.Range("F208").Value = 100112039
$ws.Range("G208").Value = "Ciboulette"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 240
$ws.Range("K208").Value = 2500
$ws.Range("L208").Value = 2500
$ws.Range("M208").Value = 2500
$ws.Range("N208").Value = "$/docena de atados"
$ws.Range("O208").Value = "Región Metropolitana"
$ws.Range("P208").Value = 833
$ws.Range("Q208").Value = 3
$ws.Range("R208").Value = "Hortaliza"
